$d = $word.ActiveDocument

# Locate the end of the previous "HOJS-022" user story (its last line is the
# "Para: ..." paragraph) and then skip one further empty paragraph (the
# spacer that carries only an rPr/sz) to land on the first of the run of
# truly empty <w:p/> paragraphs at the end of the document. That first
# empty paragraph is the one that gets replaced by the new HUJS-023 story.
$anchorText = "Para: llevar un registro de sus compras y datos."

$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*$anchorText*") {
        $target = $para.Next().Next()
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate insertion point for HUJS-023"
}
if ($target.Range.Text -ne "" -and $target.Range.Text -ne "`r") {
    throw "Expected an empty paragraph at the insertion point, found: [$($target.Range.Text)]"
}

$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$newStoryXml = @"
<w:p xmlns:w="$w"><w:r><w:rPr/><w:t>HUJS-023 Servicios de clientes</w:t></w:r></w:p><w:p xmlns:w="$w"><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">Como: Desarrollador </w:t></w:r></w:p><w:p xmlns:w="$w"><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:r><w:rPr/><w:t>Quiero: ofrecer a los clientes membrecías apropiadas para sus necesidades</w:t></w:r></w:p><w:p xmlns:w="$w"><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:r><w:rPr/><w:t>Para: poder mejorar los servicios y la retención de clientes en el concesionario.</w:t></w:r></w:p>
"@

$target.Range.InsertXML($newStoryXml)
